# Edit script: applies the changes described by the diff to before.xlsx
#  1) Swap match data (columns F:V) between row 91 and row 92
#  2) Append two new match rows (157 and 158) at the end of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the content of columns F:V between rows 91 and 92
#    (columns A:E - index/pais/torneio/temporada/data_partida - stay as-is)
# ---------------------------------------------------------------------------
$row91 = @($ws.Range("F91").Value2, $ws.Range("G91").Value2, $ws.Range("H91").Value2, $ws.Range("I91").Value2, $ws.Range("J91").Value2, $ws.Range("K91").Value2, $ws.Range("L91").Value2, $ws.Range("M91").Value2, $ws.Range("N91").Value2, $ws.Range("O91").Value2, $ws.Range("P91").Value2, $ws.Range("Q91").Value2, $ws.Range("R91").Value2, $ws.Range("S91").Value2, $ws.Range("T91").Value2, $ws.Range("U91").Value2, $ws.Range("V91").Value2)
$row92 = @($ws.Range("F92").Value2, $ws.Range("G92").Value2, $ws.Range("H92").Value2, $ws.Range("I92").Value2, $ws.Range("J92").Value2, $ws.Range("K92").Value2, $ws.Range("L92").Value2, $ws.Range("M92").Value2, $ws.Range("N92").Value2, $ws.Range("O92").Value2, $ws.Range("P92").Value2, $ws.Range("Q92").Value2, $ws.Range("R92").Value2, $ws.Range("S92").Value2, $ws.Range("T92").Value2, $ws.Range("U92").Value2, $ws.Range("V92").Value2)

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "91").Value = $row92[$i]
    $ws.Range($cols[$i] + "92").Value = $row91[$i]
}

# ---------------------------------------------------------------------------
# 2) Append two new rows (157, 158) with the same look/formatting as the
#    existing data rows (row 156 is used as the formatting template).
# ---------------------------------------------------------------------------
$ws.Range("A156:V156").Copy()
$ws.Range("A157:V158").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Positional-parameter function (named `-param value` binding is not
# reliable in this PowerShell engine, so positional args are used instead).
function Set-MatchRow($rowNum, $indice, $dataPartida, $home, $homeGols, $away, $awayGols, $homeOpenOdds, $homeOpenDh, $homeCloseOdds, $homeCloseDh, $drawOpenOdds, $drawOpenDh, $drawCloseOdds, $drawCloseDh, $awayOpenOdds, $awayOpenDh, $awayCloseOdds, $awayCloseDh, $url) {
    $ws.Range("A$rowNum").Value = $indice
    $ws.Range("B$rowNum").Value = "turkey"
    $ws.Range("C$rowNum").Value = "super-lig"
    $ws.Range("D$rowNum").Value = "2023-2024"
    $ws.Range("E$rowNum").Value = $dataPartida
    $ws.Range("F$rowNum").Value = $home
    $ws.Range("G$rowNum").Value = $homeGols
    $ws.Range("H$rowNum").Value = $away
    $ws.Range("I$rowNum").Value = $awayGols
    $ws.Range("J$rowNum").Value = $homeOpenOdds
    $ws.Range("K$rowNum").Value = $homeOpenDh
    $ws.Range("L$rowNum").Value = $homeCloseOdds
    $ws.Range("M$rowNum").Value = $homeCloseDh
    $ws.Range("N$rowNum").Value = $drawOpenOdds
    $ws.Range("O$rowNum").Value = $drawOpenDh
    $ws.Range("P$rowNum").Value = $drawCloseOdds
    $ws.Range("Q$rowNum").Value = $drawCloseDh
    $ws.Range("R$rowNum").Value = $awayOpenOdds
    $ws.Range("S$rowNum").Value = $awayOpenDh
    $ws.Range("T$rowNum").Value = $awayCloseOdds
    $ws.Range("U$rowNum").Value = $awayCloseDh
    $ws.Range("V$rowNum").Value = $url
}

Set-MatchRow 157 156 45281.625 "Rizespor" 5 "Pendikspor" 1 1.92 "14/12/2023 15:12" 1.94 "21/12/2023 14:59" 3.77 "14/12/2023 15:12" 3.78 "21/12/2023 14:59" 3.94 "14/12/2023 15:12" 4.01 "21/12/2023 14:59" "https://www.betexplorer.com/football/turkey/super-lig/rizespor-pendikspor/dU4eMYXg/"

Set-MatchRow 158 157 45281.625 "Samsunspor" 1 "Konyaspor" 1 2.12 "14/12/2023 15:12" 2.23 "21/12/2023 14:59" 3.46 "14/12/2023 15:12" 3.29 "21/12/2023 14:57" 3.57 "14/12/2023 15:12" 3.64 "21/12/2023 14:59" "https://www.betexplorer.com/football/turkey/super-lig/samsunspor-konyaspor/vgYU5DmD/"

Write-Host "Done."
